# Auto-generated: update market-price / profit figures across all 8 sheets
# (refreshed by the scheduled pricing runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 504.85715
$ws.Range("I2").Value = 592.1
$ws.Range("J2").Value = 286.75
$ws.Range("K2").Value = 592.1
$ws.Range("L2").Value = 286.75
$ws.Range("M2").Value = -479.1
$ws.Range("N2").Value = -512.75
# Row 43
$ws.Range("H43").Value = 3425
$ws.Range("I43").Value = 3425
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3425
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3356
$ws.Range("N43").ClearContents()
# Row 62
$ws.Range("H62").Value = 21951.268
$ws.Range("I62").Value = 15297.8
$ws.Range("J62").Value = 25278
$ws.Range("K62").Value = 15297.8
$ws.Range("L62").Value = 25278
$ws.Range("M62").Value = -14673.8
$ws.Range("N62").Value = -26526
# Row 65
$ws.Range("H65").Value = 21951.268
$ws.Range("I65").Value = 15297.8
$ws.Range("J65").Value = 25278
$ws.Range("K65").Value = 76489
$ws.Range("L65").Value = 126390
$ws.Range("M65").Value = -73369
$ws.Range("N65").Value = -132630
# Row 70
$ws.Range("H70").Value = 3566.5557
$ws.Range("J70").Value = 2866.6667
$ws.Range("L70").Value = 8600.000100000001
$ws.Range("N70").Value = -9140.000100000001
# Row 73
$ws.Range("H73").Value = 3566.5557
$ws.Range("J73").Value = 2866.6667
$ws.Range("L73").Value = 8600.000100000001
$ws.Range("N73").Value = -10472.0001
# Row 86
$ws.Range("H86").Value = 30976.467
$ws.Range("I86").Value = 1249.4
$ws.Range("J86").Value = 45840
$ws.Range("K86").Value = 1249.4
$ws.Range("L86").Value = 45840
$ws.Range("M86").Value = -126.4000000000001
$ws.Range("N86").Value = -48086
# Row 89
$ws.Range("H89").Value = 30976.467
$ws.Range("I89").Value = 1249.4
$ws.Range("J89").Value = 45840
$ws.Range("K89").Value = 6247
$ws.Range("L89").Value = 229200
$ws.Range("M89").Value = -631
$ws.Range("N89").Value = -240432
# Row 138
$ws.Range("H138").Value = 4502.4673
$ws.Range("I138").Value = 1517.8
$ws.Range("J138").Value = 4866.451
$ws.Range("K138").Value = 4553.4
$ws.Range("L138").Value = 14599.353
$ws.Range("M138").Value = 586.6000000000004
$ws.Range("N138").Value = -24879.353

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4711.826
$ws.Range("J61").Value = 3273.8333
$ws.Range("L61").Value = 3273.8333
$ws.Range("N61").Value = -3697.8333
# Row 63
$ws.Range("H63").Value = 1850
$ws.Range("I63").Value = 1850
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1850
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1164
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 1850
$ws.Range("I66").Value = 1850
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9250
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5818
$ws.Range("N66").ClearContents()
# Row 122
$ws.Range("H122").Value = 503709.03
$ws.Range("J122").Value = 2005425.9
$ws.Range("L122").Value = 6016277.699999999
$ws.Range("N122").Value = -6021177.699999999
# Row 136
$ws.Range("H136").Value = 4711.826
$ws.Range("J136").Value = 3273.8333
$ws.Range("L136").Value = 9821.499899999999
$ws.Range("N136").Value = -14921.4999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2289.9666
$ws.Range("I94").Value = 1845.4166
$ws.Range("K94").Value = 1845.4166
$ws.Range("M94").Value = -1394.4166

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3448.111
$ws.Range("I16").Value = 3776.1428
$ws.Range("K16").Value = 3776.1428
$ws.Range("M16").Value = -3489.1428
# Row 99
$ws.Range("H99").Value = 507469.8
$ws.Range("I99").Value = 840583
$ws.Range("J99").Value = 7800
$ws.Range("K99").Value = 840583
$ws.Range("L99").Value = 7800
$ws.Range("M99").Value = -839085
$ws.Range("N99").Value = -10796
# Row 113
$ws.Range("H113").Value = 3448.111
$ws.Range("I113").Value = 3776.1428
$ws.Range("K113").Value = 3776.1428
$ws.Range("M113").Value = -1606.1428
# Row 126
$ws.Range("H126").Value = 507469.8
$ws.Range("I126").Value = 840583
$ws.Range("J126").Value = 7800
$ws.Range("K126").Value = 2521749
$ws.Range("L126").Value = 23400
$ws.Range("M126").Value = -2519279
$ws.Range("N126").Value = -28340
# Row 132
$ws.Range("H132").Value = 22596.924
$ws.Range("I132").Value = 8255.299999999999
$ws.Range("K132").Value = 24765.9
$ws.Range("M132").Value = -22235.9

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 3675.6
$ws.Range("J113").Value = 4093.75
$ws.Range("L113").Value = 12281.25
$ws.Range("N113").Value = -16621.25
# Row 122
$ws.Range("H122").Value = 1852.5278
$ws.Range("I122").Value = 894.2222
$ws.Range("J122").Value = 2171.963
$ws.Range("K122").Value = 8047.999800000001
$ws.Range("L122").Value = 19547.667
$ws.Range("M122").Value = -5597.999800000001
$ws.Range("N122").Value = -24447.667
# Row 137
$ws.Range("H137").Value = 2645
$ws.Range("I137").Value = 1409.6154
$ws.Range("J137").Value = 7998.3335
$ws.Range("K137").Value = 4228.8462
$ws.Range("L137").Value = 23995.0005
$ws.Range("M137").Value = 871.1538
$ws.Range("N137").Value = -34195.00049999999

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 5589411.5
$ws.Range("I11").Value = 7455000
$ws.Range("J11").Value = 2169166.8
$ws.Range("K11").Value = 7455000
$ws.Range("L11").Value = 2169166.8
$ws.Range("M11").Value = -7454861
$ws.Range("N11").Value = -2169444.8
# Row 12
$ws.Range("H12").Value = 4908.1665
$ws.Range("I12").Value = 4450
$ws.Range("J12").Value = 4999.8
$ws.Range("K12").Value = 4450
$ws.Range("L12").Value = 4999.8
$ws.Range("N12").Value = -5279.8
$ws.Range("M12").Value = -4310
# Row 123
$ws.Range("H123").Value = 17312.438
$ws.Range("J123").Value = 17312.438
$ws.Range("L123").Value = 17312.438
$ws.Range("N123").Value = -22212.438
# Row 126
$ws.Range("H126").Value = 16141.667
$ws.Range("J126").Value = 9779.833000000001
$ws.Range("L126").Value = 29339.499
$ws.Range("N126").Value = -34279.499

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 15002.5
$ws.Range("J3").Value = 15002.5
$ws.Range("L3").Value = 15002.5
$ws.Range("N3").Value = -15226.5
# Row 15
$ws.Range("H15").Value = 15002.5
$ws.Range("J15").Value = 15002.5
$ws.Range("L15").Value = 15002.5
$ws.Range("N15").Value = -15342.5
# Row 68
$ws.Range("H68").Value = 4094.7368
$ws.Range("I68").Value = 3660.3
$ws.Range("J68").Value = 4577.4443
$ws.Range("K68").Value = 3660.3
$ws.Range("L68").Value = 4577.4443
$ws.Range("M68").Value = -2911.3
$ws.Range("N68").Value = -6075.4443
# Row 71
$ws.Range("H71").Value = 4094.7368
$ws.Range("I71").Value = 3660.3
$ws.Range("J71").Value = 4577.4443
$ws.Range("K71").Value = 18301.5
$ws.Range("L71").Value = 22887.2215
$ws.Range("M71").Value = -14557.5
$ws.Range("N71").Value = -30375.2215

$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 2445000
$ws.Range("I9").Value = 1676666.6
$ws.Range("K9").Value = 1676666.6
$ws.Range("M9").Value = -1676526.6
# Row 62
$ws.Range("H62").Value = 143690.9
$ws.Range("I62").Value = 346133.56
$ws.Range("J62").Value = 3538.3076
$ws.Range("K62").Value = 346133.56
$ws.Range("L62").Value = 3538.3076
$ws.Range("M62").Value = -345509.56
$ws.Range("N62").Value = -4786.3076
# Row 65
$ws.Range("H65").Value = 143690.9
$ws.Range("I65").Value = 346133.56
$ws.Range("J65").Value = 3538.3076
$ws.Range("K65").Value = 1730667.8
$ws.Range("L65").Value = 17691.538
$ws.Range("M65").Value = -1727547.8
$ws.Range("N65").Value = -23931.538
# Row 100
$ws.Range("H100").Value = 20671.852
$ws.Range("I100").Value = 11723.429
$ws.Range("K100").Value = 23446.858
$ws.Range("M100").Value = -22905.858
# Row 113
$ws.Range("H113").Value = 7056
$ws.Range("I113").Value = 4189.6
$ws.Range("J113").Value = 8847.5
$ws.Range("K113").Value = 12568.8
$ws.Range("L113").Value = 26542.5
$ws.Range("M113").Value = -10398.8
$ws.Range("N113").Value = -30882.5
# Row 122
$ws.Range("H122").Value = 28487.143
$ws.Range("J122").Value = 78856.42999999999
$ws.Range("L122").Value = 236569.29
$ws.Range("N122").Value = -241469.29

